$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 1001008
$ws.Cells.Item(6, 9).Value = 1250897.5
$ws.Cells.Item(6, 10).Value = 1450
$ws.Cells.Item(6, 11).Value = 3752692.5
$ws.Cells.Item(6, 12).Value = 4350
$ws.Cells.Item(6, 13).Value = -3752580.5
$ws.Cells.Item(6, 14).Value = -4574
$ws.Cells.Item(28, 8).Value = 42690.707
$ws.Cells.Item(28, 9).Value = 48208.43
$ws.Cells.Item(28, 11).Value = 48208.43
$ws.Cells.Item(28, 13).Value = -47723.43
$ws.Cells.Item(74, 8).Value = 6892.1816
$ws.Cells.Item(74, 9).Value = 4977.8823
$ws.Cells.Item(74, 11).Value = 4977.8823
$ws.Cells.Item(74, 13).Value = -4041.8823
$ws.Cells.Item(76, 8).Value = 7151.6665
$ws.Cells.Item(76, 9).Value = 2953
$ws.Cells.Item(76, 11).Value = 2953
$ws.Cells.Item(76, 13).Value = -2638
$ws.Cells.Item(77, 8).Value = 6892.1816
$ws.Cells.Item(77, 9).Value = 4977.8823
$ws.Cells.Item(77, 11).Value = 24889.4115
$ws.Cells.Item(77, 13).Value = -20209.4115
$ws.Cells.Item(79, 8).Value = 7151.6665
$ws.Cells.Item(79, 9).Value = 2953
$ws.Cells.Item(79, 11).Value = 2953
$ws.Cells.Item(79, 13).Value = -1861
$ws.Cells.Item(92, 8).Value = 1416.6086
$ws.Cells.Item(92, 9).Value = 662.8182
$ws.Cells.Item(92, 11).Value = 662.8182
$ws.Cells.Item(92, 13).Value = 585.1818
$ws.Cells.Item(107, 8).Value = 2058
$ws.Cells.Item(107, 9).Value = 1983.2
$ws.Cells.Item(107, 11).Value = 1983.2
$ws.Cells.Item(107, 13).Value = -63.20000000000005
$ws.Cells.Item(137, 8).Value = 2856.1365
$ws.Cells.Item(137, 10).Value = 3575.5217
$ws.Cells.Item(137, 12).Value = 10726.5651
$ws.Cells.Item(137, 14).Value = -15826.5651
$ws.Cells.Item(138, 8).Value = 3666.5894
$ws.Cells.Item(138, 10).Value = 3839.7273
$ws.Cells.Item(138, 12).Value = 11519.1819
$ws.Cells.Item(138, 14).Value = -21799.1819
$ws.Cells.Item(140, 8).Value = 69999
$ws.Cells.Item(140, 10).Value = 69999
$ws.Cells.Item(140, 12).Value = 69999
$ws.Cells.Item(140, 14).Value = -80359

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value = 33915.152
$ws.Cells.Item(43, 10).Value = 35396.11
$ws.Cells.Item(43, 12).Value = 35396.11
$ws.Cells.Item(43, 14).Value = -36022.11
$ws.Cells.Item(61, 8).Value = 7684.7144
$ws.Cells.Item(61, 9).Value = 7887
$ws.Cells.Item(61, 11).Value = 7887
$ws.Cells.Item(61, 13).Value = -7675
$ws.Cells.Item(97, 8).Value = 2375
$ws.Cells.Item(97, 9).Value = 2375
$ws.Cells.Item(97, 11).Value = 2375
$ws.Cells.Item(97, 13).Value = -1879
$ws.Cells.Item(106, 8).Value = 15000
$ws.Cells.Item(106, 9).Value = 15000
$ws.Cells.Item(106, 11).Value = 15000
$ws.Cells.Item(106, 13).Value = -13738
$ws.Cells.Item(122, 8).Value = 90912000
$ws.Cells.Item(122, 9).Value = 3196.7
$ws.Cells.Item(122, 11).Value = 9590.099999999999
$ws.Cells.Item(122, 13).Value = -7140.099999999999
$ws.Cells.Item(132, 8).Value = 5291.1304
$ws.Cells.Item(132, 9).Value = 2214.0293
$ws.Cells.Item(132, 10).Value = 14009.583
$ws.Cells.Item(132, 11).Value = 6642.0879
$ws.Cells.Item(132, 12).Value = 42028.749
$ws.Cells.Item(132, 13).Value = -4112.0879
$ws.Cells.Item(132, 14).Value = -47088.749
$ws.Cells.Item(136, 8).Value = 7684.7144
$ws.Cells.Item(136, 9).Value = 7887
$ws.Cells.Item(136, 11).Value = 23661
$ws.Cells.Item(136, 13).Value = -21111

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 979.8
$ws.Cells.Item(80, 10).Value = 966.3333
$ws.Cells.Item(80, 12).Value = 966.3333
$ws.Cells.Item(80, 14).Value = -2962.3333
$ws.Cells.Item(83, 8).Value = 979.8
$ws.Cells.Item(83, 10).Value = 966.3333
$ws.Cells.Item(83, 12).Value = 4831.6665
$ws.Cells.Item(83, 14).Value = -14815.6665
$ws.Cells.Item(94, 8).Value = 1799.5
$ws.Cells.Item(94, 9).Value = 1799.5
$ws.Cells.Item(94, 11).Value = 1799.5
$ws.Cells.Item(94, 13).Value = -1348.5
$ws.Cells.Item(107, 8).Value = 1649.6316
$ws.Cells.Item(107, 9).Value = 1352.3889
$ws.Cells.Item(107, 11).Value = 1352.3889
$ws.Cells.Item(107, 13).Value = 567.6111000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(33, 8).Value = 0
$ws.Cells.Item(33, 9).Value = 0
$ws.Cells.Item(33, 11).Value = 0
$ws.Cells.Item(33, 13).ClearContents()
$ws.Cells.Item(86, 8).Value = 11854.286
$ws.Cells.Item(86, 9).Value = 5993
$ws.Cells.Item(86, 10).Value = 19669.334
$ws.Cells.Item(86, 11).Value = 5993
$ws.Cells.Item(86, 12).Value = 19669.334
$ws.Cells.Item(86, 13).Value = -4870
$ws.Cells.Item(86, 14).Value = -21915.334
$ws.Cells.Item(89, 8).Value = 11854.286
$ws.Cells.Item(89, 9).Value = 5993
$ws.Cells.Item(89, 10).Value = 19669.334
$ws.Cells.Item(89, 11).Value = 29965
$ws.Cells.Item(89, 12).Value = 98346.67
$ws.Cells.Item(89, 13).Value = -24349
$ws.Cells.Item(89, 14).Value = -109578.67
$ws.Cells.Item(133, 8).Value = 52698
$ws.Cells.Item(133, 10).Value = 60622.625
$ws.Cells.Item(133, 12).Value = 60622.625
$ws.Cells.Item(133, 14).Value = -65682.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(102, 8).Value = 2446.75
$ws.Cells.Item(122, 8).Value = 1070143.8
$ws.Cells.Item(122, 10).Value = 2274.8235
$ws.Cells.Item(122, 12).Value = 20473.4115
$ws.Cells.Item(122, 14).Value = -25373.4115
$ws.Cells.Item(132, 8).Value = 3811.1072
$ws.Cells.Item(132, 9).Value = 2858.9285
$ws.Cells.Item(132, 10).Value = 4763.2856
$ws.Cells.Item(132, 11).Value = 25730.3565
$ws.Cells.Item(132, 12).Value = 42869.5704
$ws.Cells.Item(132, 13).Value = -23200.3565
$ws.Cells.Item(132, 14).Value = -47929.5704

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 11045.286
$ws.Cells.Item(70, 9).Value = 7566.2144
$ws.Cells.Item(70, 11).Value = 7566.2144
$ws.Cells.Item(70, 13).Value = -7296.2144
$ws.Cells.Item(73, 8).Value = 11045.286
$ws.Cells.Item(73, 9).Value = 7566.2144
$ws.Cells.Item(73, 11).Value = 7566.2144
$ws.Cells.Item(73, 13).Value = -6630.2144
$ws.Cells.Item(97, 8).Value = 1198.091
$ws.Cells.Item(97, 9).Value = 1198.091
$ws.Cells.Item(97, 11).Value = 1198.091
$ws.Cells.Item(97, 13).Value = -702.0909999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5966.5264
$ws.Cells.Item(7, 9).Value = 3726.8386
$ws.Cells.Item(7, 10).Value = 15885.143
$ws.Cells.Item(7, 11).Value = 3726.8386
$ws.Cells.Item(7, 12).Value = 15885.143
$ws.Cells.Item(7, 13).Value = -3614.8386
$ws.Cells.Item(7, 14).Value = -16109.143
$ws.Cells.Item(16, 8).Value = 0
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 13).ClearContents()
$ws.Cells.Item(82, 8).Value = 5279.4644
$ws.Cells.Item(82, 9).Value = 5039.727
$ws.Cells.Item(82, 11).Value = 5039.727
$ws.Cells.Item(82, 13).Value = -4678.727
$ws.Cells.Item(85, 8).Value = 5279.4644
$ws.Cells.Item(85, 9).Value = 5039.727
$ws.Cells.Item(85, 11).Value = 5039.727
$ws.Cells.Item(85, 13).Value = -3791.727
$ws.Cells.Item(126, 8).Value = 5966.5264
$ws.Cells.Item(126, 9).Value = 3726.8386
$ws.Cells.Item(126, 10).Value = 15885.143
$ws.Cells.Item(126, 11).Value = 11180.5158
$ws.Cells.Item(126, 12).Value = 47655.429
$ws.Cells.Item(126, 13).Value = -8710.515800000001
$ws.Cells.Item(126, 14).Value = -52595.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 4180.5454
$ws.Cells.Item(81, 9).Value = 2698.4
$ws.Cells.Item(81, 11).Value = 5396.8
$ws.Cells.Item(81, 13).Value = -4335.8
$ws.Cells.Item(84, 8).Value = 4180.5454
$ws.Cells.Item(84, 9).Value = 2698.4
$ws.Cells.Item(84, 11).Value = 26984
$ws.Cells.Item(84, 13).Value = -21680
$ws.Cells.Item(107, 8).Value = 818.8570999999999
$ws.Cells.Item(107, 9).Value = 706.9
$ws.Cells.Item(107, 11).Value = 2120.7
$ws.Cells.Item(107, 13).Value = -200.6999999999998
$ws.Cells.Item(126, 8).Value = 2770.75
$ws.Cells.Item(126, 9).Value = 2313.8
$ws.Cells.Item(126, 11).Value = 6941.400000000001
$ws.Cells.Item(126, 13).Value = -4471.400000000001
$ws.Cells.Item(136, 8).Value = 3414.84
$ws.Cells.Item(136, 9).Value = 2084.818
$ws.Cells.Item(136, 11).Value = 6254.454000000001
$ws.Cells.Item(136, 13).Value = -3704.454000000001
$ws.Cells.Item(140, 8).Value = 66628
$ws.Cells.Item(140, 9).Value = 66590
$ws.Cells.Item(140, 10).Value = 66666
$ws.Cells.Item(140, 11).Value = 66590
$ws.Cells.Item(140, 12).Value = 66666
$ws.Cells.Item(140, 13).Value = -61410
$ws.Cells.Item(140, 14).Value = -77026
$ws.Cells.Item(141, 8).Value = 103150.63
$ws.Cells.Item(141, 10).Value = 103150.63
$ws.Cells.Item(141, 12).Value = 103150.63
$ws.Cells.Item(141, 14).Value = -113510.63

Write-Output "Applied all changes"